$wb = $excel.ActiveWorkbook

# ALC row 33 (Leve Item ID 5512)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 469.6842
$ws.Range("I33").Value = 125.083336
$ws.Range("J33").Value = 1060.4286
$ws.Range("K33").Value = 125.083336
$ws.Range("L33").Value = 1060.4286
$ws.Range("M33").Value = 103.916664
$ws.Range("N33").Value = -1518.4286

# ALC row 106 (Leve Item ID 19903)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 4168778.8
$ws.Range("I106").Value = 4168778.8
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 4168778.8
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -4168147.8
$ws.Range("N106").Value = $null

# ALC row 130 (Leve Item ID 34691)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 116994
$ws.Range("J130").Value = 116994
$ws.Range("L130").Value = 116994
$ws.Range("N130").Value = -127034

# ALC row 132 (Leve Item ID 44049)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1650.16
$ws.Range("I132").Value = 1650.16
$ws.Range("K132").Value = 4950.48
$ws.Range("M132").Value = -2420.48

# ALC row 138 (Leve Item ID 44169)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2711.6611
$ws.Range("J138").Value = 3117.6287
$ws.Range("L138").Value = 9352.8861
$ws.Range("N138").Value = -19632.8861

# ARM row 44 (Leve Item ID 3861)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 68988
$ws.Range("J44").Value = 68976
$ws.Range("L44").Value = 68976
$ws.Range("N44").Value = -69952

# ARM row 55 (Leve Item ID 2830)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 61422.285
$ws.Range("J55").Value = 69984.664
$ws.Range("L55").Value = 69984.664
$ws.Range("N55").Value = -70614.664

# ARM row 61 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3910.6086
$ws.Range("I61").Value = 3876.7693
$ws.Range("J61").Value = 3954.6
$ws.Range("K61").Value = 3876.7693
$ws.Range("L61").Value = 3954.6
$ws.Range("M61").Value = -3664.7693
$ws.Range("N61").Value = -4378.6

# ARM row 63 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3833.3333
$ws.Range("I63").Value = 2375
$ws.Range("J63").Value = 5000
$ws.Range("K63").Value = 2375
$ws.Range("L63").Value = 5000
$ws.Range("M63").Value = -1689
$ws.Range("N63").Value = -6372

# ARM row 66 (Leve Item ID 12528)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H66").Value = 3833.3333
$ws.Range("I66").Value = 2375
$ws.Range("J66").Value = 5000
$ws.Range("K66").Value = 11875
$ws.Range("L66").Value = 25000
$ws.Range("M66").Value = -8443
$ws.Range("N66").Value = -31864

# ARM row 132 (Leve Item ID 43997)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3005.8635
$ws.Range("I132").Value = 2931.2942
$ws.Range("J132").Value = 3259.4
$ws.Range("K132").Value = 8793.882599999999
$ws.Range("L132").Value = 9778.200000000001
$ws.Range("M132").Value = -6263.882599999999
$ws.Range("N132").Value = -14838.2

# ARM row 136 (Leve Item ID 43999)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 3910.6086
$ws.Range("I136").Value = 3876.7693
$ws.Range("J136").Value = 3954.6
$ws.Range("K136").Value = 11630.3079
$ws.Range("L136").Value = 11863.8
$ws.Range("M136").Value = -9080.3079
$ws.Range("N136").Value = -16963.8

# ARM row 137 (Leve Item ID 43227)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = $null
$ws.Range("N137").Value = $null

# BSM row 86 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2743.2083
$ws.Range("I86").Value = 2644.8823
$ws.Range("J86").Value = 2982
$ws.Range("K86").Value = 2644.8823
$ws.Range("L86").Value = 2982
$ws.Range("M86").Value = -1521.8823
$ws.Range("N86").Value = -5228

# BSM row 89 (Leve Item ID 12526)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2743.2083
$ws.Range("I89").Value = 2644.8823
$ws.Range("J89").Value = 2982
$ws.Range("K89").Value = 13224.4115
$ws.Range("L89").Value = 14910
$ws.Range("M89").Value = -7608.411500000002
$ws.Range("N89").Value = -26142

# BSM row 99 (Leve Item ID 19943)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1932.9
$ws.Range("I99").Value = 1525.5714
$ws.Range("J99").Value = 2883.3333
$ws.Range("K99").Value = 1525.5714
$ws.Range("L99").Value = 2883.3333
$ws.Range("M99").Value = -27.57140000000004
$ws.Range("N99").Value = -5879.3333

# BSM row 105 (Leve Item ID 19947)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2287.1333
$ws.Range("I105").Value = 2440.4546
$ws.Range("J105").Value = 1865.5
$ws.Range("K105").Value = 2440.4546
$ws.Range("L105").Value = 1865.5
$ws.Range("M105").Value = -693.4546
$ws.Range("N105").Value = -5359.5

# CRP row 31 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9220.5625
$ws.Range("I31").Value = 2256.8333
$ws.Range("J31").Value = 13398.8
$ws.Range("K31").Value = 2256.8333
$ws.Range("L31").Value = 13398.8
$ws.Range("M31").Value = -1961.8333
$ws.Range("N31").Value = -13988.8

# CRP row 34 (Leve Item ID 44023)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 9220.5625
$ws.Range("I34").Value = 2256.8333
$ws.Range("J34").Value = 13398.8
$ws.Range("K34").Value = 2256.8333
$ws.Range("L34").Value = 13398.8
$ws.Range("M34").Value = -2054.8333
$ws.Range("N34").Value = -13802.8

# CRP row 132 (Leve Item ID 44019)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2591.0454
$ws.Range("I132").Value = 2255.2778
$ws.Range("J132").Value = 4102
$ws.Range("K132").Value = 6765.8334
$ws.Range("L132").Value = 12306
$ws.Range("M132").Value = -4235.8334
$ws.Range("N132").Value = -17366

# CRP row 134 (Leve Item ID 44020)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1113.8572
$ws.Range("I134").Value = 1113.8572
$ws.Range("K134").Value = 3341.5716
$ws.Range("M134").Value = -806.5715999999998

# CUL row 97 (Leve Item ID 19846)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 424.5
$ws.Range("I97").Value = 399.33334
$ws.Range("K97").Value = 1198.00002
$ws.Range("M97").Value = -702.0000199999999
$ws.Range("N97").Value = $null

# CUL row 98 (Leve Item ID 19843)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 2264.25
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

# GSM row 132 (Leve Item ID 44008)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2820.4546
$ws.Range("I132").Value = 2820.4546
$ws.Range("K132").Value = 8461.363799999999
$ws.Range("M132").Value = -5931.363799999999

# GSM row 133 (Leve Item ID 41854)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

# GSM row 135 (Leve Item ID 42006)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").Value = $null

# GSM row 137 (Leve Item ID 43226)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = $null

# LTW row 55 (Leve Item ID 5284)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 960.7353000000001
$ws.Range("I55").Value = 919.52
$ws.Range("J55").Value = 1075.2222
$ws.Range("K55").Value = 919.52
$ws.Range("L55").Value = 1075.2222
$ws.Range("M55").Value = -746.52
$ws.Range("N55").Value = -1421.2222

# LTW row 136 (Leve Item ID 44060)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 16496.857
$ws.Range("I136").Value = 7235.25
$ws.Range("J136").Value = 20201.5
$ws.Range("K136").Value = 21705.75
$ws.Range("L136").Value = 60604.5
$ws.Range("M136").Value = -19155.75
$ws.Range("N136").Value = -65704.5

# WVR row 132 (Leve Item ID 44029)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3139.7932
$ws.Range("I132").Value = 2907.5
$ws.Range("J132").Value = 3425.6924
$ws.Range("K132").Value = 8722.5
$ws.Range("L132").Value = 10277.0772
$ws.Range("M132").Value = -6192.5
$ws.Range("N132").Value = -15337.0772
